$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Course code / name update (A2 / C2) ---
# "I123A" -> "C123O"
# "Impacto Ambiental" -> "Seminario de comunicacion oral y escrita"
$ws.Range("A2").Value = "C123O"
$ws.Range("C2").Value = "Seminario de comunicacion oral y escrita"

# --- Selection moves from C10 to C12 ---
$ws.Range("C12").Select()

# --- Column C (3) widens from ~20.57 to ~40.86 characters ---
$ws.Columns.Item(3).ColumnWidth = 40

# --- Workbook window size (bookViews/workbookView) ---
$win = $wb.Windows.Item(1)
$win.Width = 28800
$win.Height = 12180
